# The "Importe" column (H) holds Argentine-formatted amounts that were
# stored as text, e.g. "22.064,00" ('.' = thousands separator, ',' =
# decimal separator). A handful of "Razon social" entries (column E) also
# happen to contain commas (e.g. "ALBIZZATTI, PABLO MARTIN Y FULINI,
# SERGIO RUBEN"). The scraper that produced this sheet ran every comma-
# bearing text value through a buggy float-formatting step: it dropped
# every "." and turned every "," into a ".", e.g.
#   "22.064,00"                                -> "22064.00"
#   "ALBIZZATTI, PABLO MARTIN Y FULINI, SERGIO" -> "ALBIZZATTI. PABLO MARTIN Y FULINI. SERGIO"
#   "...BENINCA MATIAS S.H."                    -> "...BENINCA MATIAS SH"
# Reproduce that (bug-for-bug) transformation on every cell that contains
# a comma, wherever it lives on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rowCount = $used.Rows.Count
$colCount = $used.Columns.Count
$startRow = $used.Row
$startCol = $used.Column

for ($r = 0; $r -lt $rowCount; $r++) {
    for ($c = 0; $c -lt $colCount; $c++) {
        $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
        $txt = $cell.Text
        if ($txt -ne $null -and $txt.Contains(",")) {
            $fixed = $txt.Replace(".", "").Replace(",", ".")
            # Force text storage so the numeric-looking results ("22064.00")
            # stay strings instead of being reinterpreted as numbers, then
            # drop back to the default style so no extra formatting sticks.
            $cell.NumberFormat = "@"
            $cell.Value = $fixed
            $cell.Style = "Normal"
        }
    }
}
